$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5466346740722656
$ws.Range("B1").Value = 1.07183563709259
$ws.Range("C1").Value = 5.129419326782227
$ws.Range("D1").Value = 4.078717231750488
$ws.Range("E1").Value = 0.9885925650596619
